# Generate Report for Handoff
#
# The file "918a4ba8-a289-4a33-ad17-42268040a38d.md" moved from
# "Handed back: in sync with en-US" to "Ready for handoff" and a new
# handoff was generated for each locale, so the "Latest Handoff Datetime"
# for that row is refreshed on each locale sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-01 06:30:08"

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-01 06:30:19"
